$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.471.32"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "2.647.27"
$ws.Range("E3").Value = "  +0.34%  "
$ws.Range("E4").Value = "  -0.07%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "602.53"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +2.03%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "146.64"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.55%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  -0.32%  "
$ws.Range("E9").Value = "  +1.30%  "
$ws.Range("E10").Value = "  -0.62%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.368"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +3.62%  "
$ws.Range("E12").Value = "  +0.20%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "27.39"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -0.82%  "
$ws.Range("D14").Value = "3.124.00"
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("D15").Value = "63.358.93"
$ws.Range("E15").Value = "  -0.04%  "
$ws.Range("E16").Value = "  -0.58%  "
$ws.Range("D17").Value = "2.649.16"
$ws.Range("E17").Value = "  -0.17%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "11.47"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +1.19%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "4.53"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +3.84%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "342.33"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.15%  "
$ws.Range("E21").Value = "  +2.87%  "
$ws.Range("E22").Value = "  +0.08%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "5.57"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -3.22%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "66.71"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -1.78%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "1.69"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.57%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "8.85"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +4.11%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "1.53"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -2.51%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "548.28"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -1.05%  "
$ws.Range("E29").Value = "  -1.87%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.34%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "7.88"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.83%  "
$ws.Range("E32").Value = "  +3.79%  "
$ws.Range("E33").Value = "  -2.26%  "
$ws.Range("D34").Value = "0.0₃0807"
$ws.Range("E34").Value = "  -0.39%  "
$ws.Range("E35").Value = "  +8.24%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "167.29"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -4.50%  "
$ws.Range("E37").Value = "  +1.08%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.05%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "19.10"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.17%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "1.90"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +6.99%  "
$ws.Range("E41").Value = "  -0.10%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "169.19"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.83%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "3.77"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.51%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "22.50"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +2.49%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.0576"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +4.31%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.625"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.56%  "
$ws.Range("E47").Value = "  +2.37%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.0963"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +0.56%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "18.79"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -0.29%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.78"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +4.46%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "11.28"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.62%  "
